# Update results w/ openjml handling fixed

$wb = $excel.ActiveWorkbook

# --- Sheet "all_tools": single value fix ---
$wsAll = $wb.Worksheets.Item("all_tools")
$wsAll.Range("F9").Value = 93

# --- Sheet "openjml": correlation analysis refresh for rows 9-12 ---
$wsOpenjml = $wb.Worksheets.Item("openjml")

# Row 9
$wsOpenjml.Range("F9").Value = 49
$wsOpenjml.Range("H9").Value = 100
$wsOpenjml.Range("I9").Value = -0.1777071888309779
$wsOpenjml.Range("J9").Value = 0.02109743444634812
$wsOpenjml.Range("K9").Value = -0.2338626244994551
$wsOpenjml.Range("L9").Value = 0.01919115836424928

# Row 10
$wsOpenjml.Range("F10").Value = 40
$wsOpenjml.Range("H10").Value = 50
$wsOpenjml.Range("I10").Value = -0.2254320318923924
$wsOpenjml.Range("J10").Value = 0.03750206210014283
$wsOpenjml.Range("K10").Value = -0.2757884146130811
$wsOpenjml.Range("L10").Value = 0.0525532974249969

# Row 11
$wsOpenjml.Range("F11").Value = 40
$wsOpenjml.Range("H11").Value = 50
$wsOpenjml.Range("I11").Value = -0.1794340300657417
$wsOpenjml.Range("J11").Value = 0.08321105252268295
$wsOpenjml.Range("K11").Value = -0.2382169979796553
$wsOpenjml.Range("L11").Value = 0.09572883636182701

# Row 12
$wsOpenjml.Range("F12").Value = 40
$wsOpenjml.Range("H12").Value = 50
$wsOpenjml.Range("I12").Value = 0.1869191683862415
$wsOpenjml.Range("J12").Value = 0.06934431565042681
$wsOpenjml.Range("K12").Value = 0.2610882725528407
$wsOpenjml.Range("L12").Value = 0.06704286140152567
